$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $value) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextCell 66 1 "-538"
Set-TextCell 66 2 "7/31/2025"
Set-TextCell 66 3 "Malabia 964"
Set-TextCell 66 4 "15"
Set-TextCell 66 5 "808609237"
Set-TextCell 66 6 "NEW"
Set-TextCell 66 7 "Pendiente"
Set-TextCell 66 8 "Cambiar poste mal estado por PRFV"

$ws.Cells.Item(66, 9).Value = 1

Set-TextCell 66 10 "Cambio"
Set-TextCell 66 11 "Sin equipos"
Set-TextCell 66 12 "Poste"

$ws.Cells.Item(66, 13).Value = -58.433634
$ws.Cells.Item(66, 14).Value = -34.595018

Set-TextCell 66 15 "Palermo"
Set-TextCell 66 16 "Capital Sur"
